$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 124, shifting existing rows 124:230 down to 125:231
$ws.Rows.Item(124).Insert()

# Populate the newly inserted row 124 with the new weekly price record
$ws.Cells.Item(124, 1).Value = 4
$ws.Cells.Item(124, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(124, 3).Value = "Los Lagos"
$ws.Cells.Item(124, 4).Value = 44658
$ws.Cells.Item(124, 5).Value = 10
$ws.Cells.Item(124, 6).Value = "Fruta"
$ws.Cells.Item(124, 7).Value = 100101
$ws.Cells.Item(124, 8).Value = "Berries"
$ws.Cells.Item(124, 9).Value = 100112025
$ws.Cells.Item(124, 10).Value = "Frutilla"
$ws.Cells.Item(124, 11).Value = "Sin especificar"
$ws.Cells.Item(124, 12).Value = "Primera"
$ws.Cells.Item(124, 13).Value = 500
$ws.Cells.Item(124, 14).Value = 8000
$ws.Cells.Item(124, 15).Value = 8500
$ws.Cells.Item(124, 16).Value = 8250
$ws.Cells.Item(124, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(124, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(124, 19).Value = 1179
$ws.Cells.Item(124, 20).Value = 7
